$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 (the "2006年" data row); all subsequent rows shift up by one.
$ws.Rows.Item(2).Delete()
